# Corrected comparison of strings to .equals
# Add a second worksheet ("Sheet2") right after "Sheet1", put the text
# "TEST" into its A1 cell, and make it the active/selected sheet
# (mirrors what Excel does when you insert a new sheet and leave it
# selected: tabSelected moves from Sheet1 to Sheet2 and the workbook's
# activeTab becomes index 1).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert the new worksheet immediately after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "TEST"

# Leave Sheet2 as the active sheet/tab
$ws2.Select()
